$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Bring in row 7's exact formatting (border etc.) for the new row 8 by
# copy/paste-special of formats, then fill in the new row's values.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)

$ws.Range("A8").Value = "Sanity suite"
$ws.Range("B8").Value = "Sanity tests"
$ws.Range("C8").Value = "Y"

# The new "Sanity" cells got their fill explicitly set to "No Fill" in the
# authored workbook (distinct style from the rest of the data rows).
$ws.Range("B8:C8").Interior.ColorIndex = -4142

# Selection moved on to B13 after the edit.
$ws.Range("B13").Select()
